{"js": "// \"Ch\u1ec9nh l\u1ea1i m\u1eabu 26\": remove the placeholder merge-field text\n// \"vnpt.SiteAddress\" that followed the \"\u0110\u1ecba ch\u1ec9: \" label, leaving just\n// the label itself in the paragraph.\nconst results = context.document.body.search(\"vnpt.SiteAddress\", { matchCase: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (const result of results.items) {\n  result.delete();\n}\nawait context.sync();\n", "ps1": "# \"Ch\u1ec9nh l\u1ea1i m\u1eabu 26\": remove the placeholder merge-field text\n# \"vnpt.SiteAddress\" that followed the \"\u0110\u1ecba ch\u1ec9: \" label, leaving just\n# the label itself in the paragraph.\n$d = $word.ActiveDocument\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"vnpt.SiteAddress\"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $range.Delete()\n}\n"}
